# Adapt column header formatting to respective input file names (#7)
# - rename "<col>_old" headers to "<col>_FV2310"
# - rename "<col>_new" headers to "<col>_FV2404"
# - wrap the data range in an Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) ---------------------------------------

$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

# columns A..J (1..10) hold the "_old" -> "_FV2310" headers
for ($i = 0; $i -lt $fv2310Headers.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $fv2310Headers[$i]
}

# column K (11) holds "diff" and stays untouched

# columns L..U (12..21) hold the "_new" -> "_FV2404" headers
for ($i = 0; $i -lt $fv2404Headers.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $fv2404Headers[$i]
}

# --- 2. Turn the data range into an Excel Table --------------------------

$tableRange = $ws.Range("A1:U53")
$lo = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row --------------------------------------------

$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
